# Doing Updates for Financials
# Updates the NMR yearly financials figures in place (Income Statement,
# Balance Sheet, Cash Flow Statement sections) to the latest reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NMR")

function Set-Row {
    param(
        [int]$Row,
        $D = $null,
        $E = $null,
        $F = $null,
        $G = $null,
        $H = $null,
        $I = $null,
        $J = $null
    )

    if ($null -ne $D) { $ws.Range("D$Row").Value = $D }
    if ($null -ne $E) { $ws.Range("E$Row").Value = $E }
    if ($null -ne $F) { $ws.Range("F$Row").Value = $F }
    if ($null -ne $G) { $ws.Range("G$Row").Value = $G }
    if ($null -ne $H) { $ws.Range("H$Row").Value = $H }
    if ($null -ne $I) { $ws.Range("I$Row").Value = $I }
    if ($null -ne $J) { $ws.Range("J$Row").Value = $J }
}

# ---- Income Statement ----
Set-Row 8  17828300 15508300 15576800 17452500 16559900 18802700 16739900
Set-Row 9  4295700  2823400  2959800  2950800  2484000  2407500  2855700
Set-Row 10 13532600 12684900 12617000 14501800 14075900 16395200 13884200
Set-Row 14 $null $null $null $null $null $null 4485900
Set-Row 17 14861800 12590200 14083800 14317800 13290900 16653600 15971900
Set-Row 18 2966500  2918100  1493000  3134700  3269000  2149100  768000
Set-Row 21 3614300  3560000  2211500  3848600  3988200  2977100  1678200
Set-Row 23 2966500  2918100  1493000  3134700  3269000  2149100  768000
Set-Row 24 938900   725300   204300   1091900  1312300  1193600  532500
Set-Row 26 2027600  2192800  1288800  2042900  1956700  955400   235500
Set-Row 27 1982900  2166100  1189200  2032100  1930900  969400   104700
Set-Row 33 1982900  2166100  1189200  2032100  1930900  969400   104700
Set-Row 35 1982900  2166100  1189200  2032100  1930900  969400   104700

# ---- Balance Sheet ----
Set-Row 41 24137600  24811500  33203000  14857800  16755400  12502400  15584800
Set-Row 43 12784800  11056300  12339300  13454300  11242600  9540100   77614600
Set-Row 45 2612200   2056200   2042600   4095500   3036000   2438500   "NA"
Set-Row 47 318366000 341995000 315901000 337419000 354954000 309204000 204682000
Set-Row 48 3064400   3161300   3213800   3625700   3696600   3871300   9455400
Set-Row 49 899100    947600    999200    1116300   1040900
Set-Row 52 145900    197300    326600    178300    199000    "NA"
Set-Row 54 364709000 387383000 371455000 377721000 393424000 343000000 322704000
Set-Row 58 164331000 184406000 171898000 150008000 162845000 148378000 124016000
Set-Row 59 28927900  31139900  25009400  25868000  21585800  12128700  13855800
Set-Row 61 66737900  65046500  73491200  75360100  74372600  68635000  76883800
Set-Row 62 $null "NA" "NA" "NA" 314000
Set-Row 66 339855000 362162000 347045000 353242000 370700000 322259000 303654000
Set-Row 72 15339900  15035600  13709900  12999000  11634500  10274200  9572900
Set-Row 76 24853900  25220800  24410200  24478300  22723700  20741100  19049500

# ---- Cash Flow Statement ----
Set-Row 81  1982900  2166100   1189200  2032100   1930900  969400   104700
Set-Row 83  647100   641200    717700   713100    718400   827100   909200
Set-Row 89  -4029100 11797400  11194900 -696300   4135100  4967500  2629400
Set-Row 91  -2577900 -2828400  -2935500 -1893600  -1937600 -2458700 -1650400
Set-Row 94  -507800  -1067200  -214300  111500    -932900  -1450800 89900
Set-Row 96  -634600  -387200   -748400  -500100   -469600  -133200  -262800
Set-Row 100 3373400  -19261000 8916900  -1611000  2616000  -6342700 -7632600
Set-Row 101 -483600  38400     -363400  619400    371400   426500   -57100
Set-Row 102 -1647100 -8492400  19534100 -1576400  6189700  -2399500 -4970400
